$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 213.93333
$ws.Range("I9").Value = 88.416664
$ws.Range("K9").Value = 88.416664
$ws.Range("M9").Value = 80.583336
$ws.Range("H62").Value = 7263.1816
$ws.Range("I62").Value = 8863.125
$ws.Range("K62").Value = 8863.125
$ws.Range("M62").Value = -8239.125
$ws.Range("H65").Value = 7263.1816
$ws.Range("I65").Value = 8863.125
$ws.Range("K65").Value = 44315.625
$ws.Range("M65").Value = -41195.625
$ws.Range("H76").Value = 3850.4644
$ws.Range("I76").Value = 2988.3333
$ws.Range("K76").Value = 2988.3333
$ws.Range("M76").Value = -2673.3333
$ws.Range("H79").Value = 3850.4644
$ws.Range("I79").Value = 2988.3333
$ws.Range("K79").Value = 2988.3333
$ws.Range("M79").Value = -1896.3333
$ws.Range("H118").Value = 219.75
$ws.Range("I118").Value = 193.33333
$ws.Range("K118").Value = 579.99999
$ws.Range("M118").Value = 1077.00001
$ws.Range("H137").Value = 21281864
$ws.Range("I137").Value = 62501736
$ws.Range("J137").Value = 7092.2583
$ws.Range("K137").Value = 187505208
$ws.Range("L137").Value = 21276.7749
$ws.Range("M137").Value = -187502658
$ws.Range("N137").Value = -26376.7749
$ws.Range("H138").Value = 2583.5334
$ws.Range("J138").Value = 2593.0908
$ws.Range("L138").Value = 7779.2724
$ws.Range("N138").Value = -18059.2724

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 846440.1
$ws.Range("I74").Value = 1070562.4
$ws.Range("J74").Value = 13986.071
$ws.Range("K74").Value = 1070562.4
$ws.Range("L74").Value = 13986.071
$ws.Range("M74").Value = -1069688.4
$ws.Range("N74").Value = -15734.071
$ws.Range("H77").Value = 846440.1
$ws.Range("I77").Value = 1070562.4
$ws.Range("J77").Value = 13986.071
$ws.Range("K77").Value = 5352812
$ws.Range("L77").Value = 69930.355
$ws.Range("M77").Value = -5348444
$ws.Range("N77").Value = -78666.355
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""
$ws.Range("H109").Value = 63125.668
$ws.Range("J109").Value = 63125.668
$ws.Range("L109").Value = 63125.668
$ws.Range("N109").Value = -65899.66800000001
$ws.Range("H112").Value = 25844
$ws.Range("J112").Value = 25844
$ws.Range("L112").Value = 25844
$ws.Range("N112").Value = -28798
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""
$ws.Range("H122").Value = 1167.7778
$ws.Range("I122").Value = 876.25
$ws.Range("K122").Value = 2628.75
$ws.Range("M122").Value = -178.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 252050.75
$ws.Range("J16").Value = 252373.75
$ws.Range("L16").Value = 252373.75
$ws.Range("N16").Value = -252947.75
$ws.Range("H31").Value = 2648319.2
$ws.Range("I31").Value = 2648319.2
$ws.Range("K31").Value = 2648319.2
$ws.Range("M31").Value = -2648024.2
$ws.Range("H34").Value = 2648319.2
$ws.Range("I34").Value = 2648319.2
$ws.Range("K34").Value = 2648319.2
$ws.Range("M34").Value = -2648117.2
$ws.Range("H59").Value = 45333
$ws.Range("J59").Value = 45333
$ws.Range("L59").Value = 45333
$ws.Range("N59").Value = -47623
$ws.Range("H60").Value = 27500
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 27500
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 27500
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = -28522
$ws.Range("H113").Value = 252050.75
$ws.Range("J113").Value = 252373.75
$ws.Range("L113").Value = 252373.75
$ws.Range("N113").Value = -256713.75
$ws.Range("H132").Value = 1919.6471
$ws.Range("I132").Value = 1883.6
$ws.Range("J132").Value = 1971.1428
$ws.Range("K132").Value = 5650.799999999999
$ws.Range("L132").Value = 5913.428400000001
$ws.Range("M132").Value = -3120.799999999999
$ws.Range("N132").Value = -10973.4284
$ws.Range("H134").Value = 2069.2812
$ws.Range("I134").Value = 1619.1852
$ws.Range("K134").Value = 4857.5556
$ws.Range("M134").Value = -2322.5556

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7499.222
$ws.Range("I56").Value = 7499.222
$ws.Range("K56").Value = 7499.222
$ws.Range("M56").Value = -6969.222
$ws.Range("H122").Value = 949858.6
$ws.Range("I122").Value = 1613853.2
$ws.Range("J122").Value = 1295
$ws.Range("K122").Value = 14524678.8
$ws.Range("L122").Value = 11655
$ws.Range("M122").Value = -14522228.8
$ws.Range("N122").Value = -16555

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 57581
$ws.Range("J51").Value = 57581
$ws.Range("L51").Value = 57581
$ws.Range("N51").Value = -58599
$ws.Range("H70").Value = 38398
$ws.Range("I70").Value = 33828.285
$ws.Range("J70").Value = 44795.6
$ws.Range("K70").Value = 33828.285
$ws.Range("L70").Value = 44795.6
$ws.Range("M70").Value = -33558.285
$ws.Range("N70").Value = -45335.6
$ws.Range("H73").Value = 38398
$ws.Range("I73").Value = 33828.285
$ws.Range("J73").Value = 44795.6
$ws.Range("K73").Value = 33828.285
$ws.Range("L73").Value = 44795.6
$ws.Range("M73").Value = -32892.285
$ws.Range("N73").Value = -46667.6
$ws.Range("H103").Value = 74960.2
$ws.Range("J103").Value = 74960.2
$ws.Range("L103").Value = 74960.2
$ws.Range("N103").Value = -77304.2
$ws.Range("H111").Value = 43323.25
$ws.Range("J111").Value = 43323.25
$ws.Range("L111").Value = 43323.25
$ws.Range("N111").Value = -49457.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2051.1333
$ws.Range("I22").Value = 1897.6666
$ws.Range("K22").Value = 1897.6666
$ws.Range("M22").Value = -1602.6666
$ws.Range("H27").Value = 2051.1333
$ws.Range("I27").Value = 1897.6666
$ws.Range("K27").Value = 1897.6666
$ws.Range("M27").Value = -1790.6666
$ws.Range("H40").Value = 4134.4736
$ws.Range("I40").Value = 3524.375
$ws.Range("J40").Value = 7388.3335
$ws.Range("K40").Value = 3524.375
$ws.Range("L40").Value = 7388.3335
$ws.Range("M40").Value = -3388.375
$ws.Range("N40").Value = -7660.3335
$ws.Range("H46").Value = 3058.7646
$ws.Range("I46").Value = 611.1111
$ws.Range("J46").Value = 5812.375
$ws.Range("K46").Value = 611.1111
$ws.Range("L46").Value = 5812.375
$ws.Range("M46").Value = -423.1111
$ws.Range("N46").Value = -6188.375
$ws.Range("H93").Value = 1802.5161
$ws.Range("I93").Value = 1493.16
$ws.Range("J93").Value = 3091.5
$ws.Range("K93").Value = 1493.16
$ws.Range("L93").Value = 3091.5
$ws.Range("M93").Value = -245.1600000000001
$ws.Range("N93").Value = -5587.5
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = ""
$ws.Range("H110").Value = 31350
$ws.Range("J110").Value = 31350
$ws.Range("L110").Value = 31350
$ws.Range("N110").Value = -39530
$ws.Range("H132").Value = 1968978.8
$ws.Range("I132").Value = 2478436.2
$ws.Range("K132").Value = 7435308.600000001
$ws.Range("M132").Value = -7432778.600000001
$ws.Range("H136").Value = 7360151.5
$ws.Range("I136").Value = 7817188
$ws.Range("J136").Value = 6953896.5
$ws.Range("K136").Value = 23451564
$ws.Range("L136").Value = 20861689.5
$ws.Range("M136").Value = -23449014
$ws.Range("N136").Value = -20866789.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3547871
$ws.Range("I132").Value = 3969783
$ws.Range("K132").Value = 11909349
$ws.Range("M132").Value = -11906819
$ws.Range("H136").Value = 5411900.5
$ws.Range("I136").Value = 1359793.2
$ws.Range("K136").Value = 4079379.6
$ws.Range("M136").Value = -4076829.6
